$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "72.855.75"
$ws.Range("E2").Value = "  +1.81%  "
Set-TextValue $ws.Range("D3") "3.986.31"
$ws.Range("E3").Value = "  +0.01%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.11%  "
Set-TextValue $ws.Range("D5") "619.84"
$ws.Range("E5").Value = "  +16.08%  "
Set-TextValue $ws.Range("D6") "168.62"
$ws.Range("E6").Value = "  +12.41%  "
Set-TextValue $ws.Range("D7") "0.684"
$ws.Range("E7").Value = "  -0.26%  "
Set-TextValue $ws.Range("D8") "0.999"
$ws.Range("E8").Value = "  -0.10%  "
Set-TextValue $ws.Range("D9") "0.759"
$ws.Range("E9").Value = "  +1.59%  "
Set-TextValue $ws.Range("D10") "0.167"
$ws.Range("E10").Value = "  -0.17%  "
Set-TextValue $ws.Range("D11") "56.07"
$ws.Range("E11").Value = "  +2.27%  "
Set-TextValue $ws.Range("D12") "0.0000316"
$ws.Range("E12").Value = "  -1.37%  "
Set-TextValue $ws.Range("D13") "11.25"
$ws.Range("E13").Value = "  +5.28%  "
Set-TextValue $ws.Range("D14") "4.621.11"
$ws.Range("E14").Value = "  +0.21%  "
Set-TextValue $ws.Range("D15") "3.964.76"
$ws.Range("E15").Value = "  -0.23%  "
Set-TextValue $ws.Range("D16") "1.28"
$ws.Range("E16").Value = "  +8.29%  "
Set-TextValue $ws.Range("D17") "14.22"
$ws.Range("E17").Value = "  +1.53%  "
Set-TextValue $ws.Range("D18") "20.69"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("E19").Value = "  +0.45%  "
Set-TextValue $ws.Range("D20") "72.580.48"
$ws.Range("E20").Value = "  +1.59%  "
Set-TextValue $ws.Range("D21") "438.33"
$ws.Range("E21").Value = "  +2.01%  "
Set-TextValue $ws.Range("D22") "4.95"
$ws.Range("E22").Value = "  +18.07%  "
Set-TextValue $ws.Range("D23") "96.32"
$ws.Range("E23").Value = "  -1.17%  "
Set-TextValue $ws.Range("D24") "3.43"
$ws.Range("E24").Value = "  -3.82%  "
Set-TextValue $ws.Range("D25") "14.54"
$ws.Range("E25").Value = "  -0.62%  "
Set-TextValue $ws.Range("D26") "4.34"
$ws.Range("E26").Value = "  +5.51%  "
Set-TextValue $ws.Range("D27") "11.40"
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("E28").Value = "  -1.74%  "
Set-TextValue $ws.Range("D29") "5.95"
$ws.Range("E29").Value = "  +0.76%  "
Set-TextValue $ws.Range("D30") "36.21"
$ws.Range("E30").Value = "  -1.12%  "
Set-TextValue $ws.Range("D31") "7.79"
$ws.Range("E31").Value = "  -4.97%  "
Set-TextValue $ws.Range("D32") "13.92"
$ws.Range("E32").Value = "  +3.56%  "
Set-TextValue $ws.Range("D33") "0.131"
$ws.Range("E33").Value = "  -4.77%  "
Set-TextValue $ws.Range("D34") "72.46"
$ws.Range("E34").Value = "  +10.33%  "
Set-TextValue $ws.Range("D35") "48.33"
$ws.Range("E35").Value = "  -7.39%  "
Set-TextValue $ws.Range("D36") "642.92"
$ws.Range("E36").Value = "  -6.23%  "
Set-TextValue $ws.Range("D37") "0.0₃0888"
$ws.Range("E37").Value = "  +8.87%  "
Set-TextValue $ws.Range("D38") "0.440"
$ws.Range("E38").Value = "  -4.89%  "
Set-TextValue $ws.Range("D40") "3.40"
$ws.Range("E40").Value = "  +6.08%  "
Set-TextValue $ws.Range("D42") "0.999"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  +0.35%  "
Set-TextValue $ws.Range("D44") "0.0489"
$ws.Range("E44").Value = "  -0.90%  "
Set-TextValue $ws.Range("D45") "10.66"
$ws.Range("E45").Value = "  +0.98%  "
Set-TextValue $ws.Range("D46") "0.150"
$ws.Range("E46").Value = "  +0.34%  "
Set-TextValue $ws.Range("D47") "2.67"
$ws.Range("E47").Value = "  -0.26%  "
Set-TextValue $ws.Range("D48") "3.45"
$ws.Range("E48").Value = "  +2.20%  "
Set-TextValue $ws.Range("D49") "3.10"
$ws.Range("E49").Value = "  +1.15%  "
Set-TextValue $ws.Range("D50") "2.908.46"
$ws.Range("E50").Value = "  +9.14%  "
Set-TextValue $ws.Range("D51") "3.42"
$ws.Range("E51").Value = "  +1.26%  "

# Row 39 and 41 swap coin identity (Kaspa <-> ThetaToken) with new price/volume data
$ws.Range("B39").Value = "ThetaToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D39") "3.42"
$ws.Range("E39").Value = "  +1.91%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D41") "0.148"
$ws.Range("E41").Value = "  -1.27%  "
